{"js": "// Fix Race Collapse Variable (race_ethn_5cat) \u2014 \"Non-Hispanic Asian\" row was\n// erroneously omitted (all zeros), and the \"Other Race (including\n// multiracial)\" row absorbed those subjects. Restore the correct N and\n// Statistic (95% CI) values for both rows in the comparative table.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No tables found in document body.\");\n}\n\nconst table = tables.items[0];\n\n// Load the full grid of cell text so we can locate the target rows by their\n// \"Category\" label (column index 1) instead of relying on a hard-coded row\n// number.\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst rowCount = rows.items.length;\nconst cells2d = [];\nfor (let r = 0; r < rowCount; r++) {\n  const rowCells = rows.items[r].cells;\n  rowCells.load(\"items\");\n  cells2d.push(rowCells);\n}\nawait context.sync();\n\nfor (let r = 0; r < rowCount; r++) {\n  for (const c of cells2d[r].items) {\n    c.load(\"value\");\n  }\n}\nawait context.sync();\n\nfunction findRowByCategory(label) {\n  for (let r = 0; r < rowCount; r++) {\n    const items = cells2d[r].items;\n    if (items.length > 1 && items[1].value.trim() === label) {\n      return r;\n    }\n  }\n  return -1;\n}\n\nconst asianRow = findRowByCategory(\"Non-Hispanic Asian\");\nconst otherRow = findRowByCategory(\"Other Race (including multiracial)\");\n\nif (asianRow === -1 || otherRow === -1) {\n  throw new Error(\"Could not locate target table rows by Category label.\");\n}\n\n// column 2: N - TRUE, column 3: Statistic (95% CI) - TRUE,\n// column 4: N - FALSE, column 5: Statistic (95% CI) - FALSE\nconst updates = [\n  { row: asianRow, col: 2, value: \"87\" },\n  { row: asianRow, col: 3, value: \"11.3 (9.3 - 13.8)\" },\n  { row: asianRow, col: 4, value: \"37\" },\n  { row: asianRow, col: 5, value: \"9.2 (6.7 - 12.5)\" },\n  { row: otherRow, col: 2, value: \"36\" },\n  { row: otherRow, col: 3, value: \"4.7 (3.4 - 6.4)\" },\n  { row: otherRow, col: 4, value: \"12\" },\n  { row: otherRow, col: 5, value: \"3.0 (1.7 - 5.2)\" },\n];\n\nfor (const u of updates) {\n  table.getCell(u.row, u.col).value = u.value;\n}\n\nawait context.sync();\n", "ps1": "# Fix Race Collapse Variable (race_ethn_5cat) - \"Non-Hispanic Asian\" row was\n# erroneously omitted (all zeros), and the \"Other Race (including\n# multiracial)\" row absorbed those subjects. Restore the correct N and\n# Statistic (95% CI) values for both rows in the comparative table.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Locate the two target rows by their \"Category\" label (column 2) instead of\n# relying on a hard-coded row index, so the script is resilient to layout\n# changes elsewhere in the table.\n$asianRow = 0\n$otherRow = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    $label = $t.Cell($r, 2).Range.Text\n    $label = $label.TrimEnd([char]13, [char]7)\n    if ($label -eq \"Non-Hispanic Asian\") {\n        $asianRow = $r\n    }\n    elseif ($label -eq \"Other Race (including multiracial)\") {\n        $otherRow = $r\n    }\n}\n\nif ($asianRow -eq 0 -or $otherRow -eq 0) {\n    throw \"Could not locate target table rows by Category label.\"\n}\n\n# Column 3: N - TRUE, Column 4: Statistic (95% CI) - TRUE,\n# Column 5: N - FALSE, Column 6: Statistic (95% CI) - FALSE\n$t.Cell($asianRow, 3).Range.Text = \"87\"\n$t.Cell($asianRow, 4).Range.Text = \"11.3 (9.3 - 13.8)\"\n$t.Cell($asianRow, 5).Range.Text = \"37\"\n$t.Cell($asianRow, 6).Range.Text = \"9.2 (6.7 - 12.5)\"\n\n$t.Cell($otherRow, 3).Range.Text = \"36\"\n$t.Cell($otherRow, 4).Range.Text = \"4.7 (3.4 - 6.4)\"\n$t.Cell($otherRow, 5).Range.Text = \"12\"\n$t.Cell($otherRow, 6).Range.Text = \"3.0 (1.7 - 5.2)\"\n"}
